$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendar2021")
$ws.Activate()

# Update cell F3: HW1 link changed to canvas quizzes URL
$ws.Range("F3").Value = "HW 1 (https://canvas.jmu.edu/courses/1775272/quizzes); PA 1(PAs/PA.php?paNumber=1)"

# Update cell D4: Informed Search description + slides/video link added
$ws.Range("D4").Value = "Informed Search = A* and Graph Search;slides(slides/03_03_InformedSearch.pdf)   video(https://canvas.jmu.edu/courses/1775272/modules)"

# Update cell F4: HW 2 link changed to canvas quizzes URL
$ws.Range("F4").Value = "HW 2 (https://canvas.jmu.edu/courses/1775272/quizzes)"

# Update cell G5: remove "PA 1;" prefix
$ws.Range("G5").Value = "HW 2;Quiz 0"

# Update the active selection on sheet1 to D5
$ws.Range("D5").Select()

# Update workbook window position
$excel.ActiveWindow.Left = 8560
$excel.ActiveWindow.Top = 4720
